$wb = $excel.ActiveWorkbook

# Add the new worksheet, placed after "Stromverbrauch"
$wsStrom = $wb.Worksheets.Item("Stromverbrauch")
$ws = $wb.Worksheets.Add([System.Type]::Missing, $wsStrom)
$ws.Name = "Spannungsregler"

# Values are entered in the exact sequence the original author used, so that
# the resulting shared-string table indices line up with the target file.

# 1) B1 = "4s zu 12.5V"
$ws.Range("B1").Value = "4s zu 12.5V"

# 2) A2 = "Int.referenz"
$ws.Range("A2").Value = "Int.referenz"
$ws.Range("B2").Value = 1.25
$ws.Range("C2").Value = 1.25
$ws.Range("D2").Value = 1.25
$ws.Range("E2").Value = "V"

# 3) C1 = "4s zu 10V"
$ws.Range("C1").Value = "4s zu 10V"

# 4) D1 = "10V zu 3V3"
$ws.Range("D1").Value = "10V zu 3V3"

# 5) A4 = "R1 (ausgewählt)"
$ws.Range("A4").Value = "R1 (ausgewählt)"
$ws.Range("B4").Value = 5100
$ws.Range("C4").Value = 5100
$ws.Range("D4").Value = 5100
$ws.Range("E4").Value = "Ohm"

# 6) A5 = "R2 (Berechnet)"
$ws.Range("A5").Value = "R2 (Berechnet)"
$ws.Range("B5").Formula = "=((B3/B2)-1)*B4"
$ws.Range("C5").Formula = "=((C3/C2)-1)*C4"
$ws.Range("D5").Formula = "=((D3/D2)-1)*D4"
$ws.Range("E5").Value = "Ohm"

# 7) A6 = "R2 (ausgewählt"
$ws.Range("A6").Value = "R2 (ausgewählt"
$ws.Range("B6").Value = 47000
$ws.Range("C6").Value = 33000
$ws.Range("D6").Value = 8200
$ws.Range("E6").Value = "Ohm"

# 8) A3 = "Ausgangsspannung Ziel"
$ws.Range("A3").Value = "Ausgangsspannung Ziel"
$ws.Range("B3").Value = 12.5
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 3.3
$ws.Range("E3").Value = "V"

# 9) A7 = "Ausgangsspannung"
$ws.Range("A7").Value = "Ausgangsspannung"
$ws.Range("B7").Formula = "=1.25*(1+(B6/B4))"
$ws.Range("C7").Formula = "=1.25*(1+(C6/C4))"
$ws.Range("D7").Formula = "=1.25*(1+(D6/D4))"
$ws.Range("E7").Value = "V"

$ws.Range("B7:D7").NumberFormat = "0.00"

# Column widths similar to author's workbook (best-fit widths from the
# original file, expressed as ColumnWidth = storedWidth - 5/6)
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668
$ws.Columns.Item(2).ColumnWidth = 9.736979166666666
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 9.592447916666666

$ws.Range("E8").Select()

$ws.Activate()
